$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Authors column (E) for these rows contained values exported with an
# inconsistent amount of whitespace after each separator comma. Clean the
# data by normalizing to two additional spaces after every comma (matches
# the rest of the cleaned dataset).
$rows = @(2,3,4,5,6,7,8,9,10,11,14,15,16,17,18,19)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 5)
    $old = $cell.Value2
    $new = $old -replace ",( +)", ',$1  '
    $cell.Value2 = $new
}
